$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the "+1" callout textbox ("TextBox 1059") that sits next to the
# PDP stream diagram. There are two shapes on the slide whose text reads
# "+1"; the one we want is the free-standing textbox at
# Left=686.6334pt, Top=477.5959pt (the other is an unrelated rectangle
# elsewhere on the slide), so disambiguate on position as well as text.
# Note: use .Equals() for the text check - PowerShell's -eq operator
# coerces numeric-looking strings (so "1" -eq "+1" is $true), which
# would also match the unrelated rectangle's plain "1" label.
$target = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text.Equals("+1") -and [math]::Round($shp.Top, 0) -eq 478) {
                $target = $shp
                break
            }
        }
    }
}

$tf = $target.TextFrame
$tr = $tf.TextRange

# Update the label text from "+1" to "+n" (the stream can now advance
# by an arbitrary n, not just 1).
$tr.Text = "+n"

# Narrow the (auto-fit) textbox to match the new, slightly narrower,
# glyph "n" in place of "1".
$target.Width = 28.172440944881888
